$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.667.76"
$ws.Range("E2").Value = "  -7.29%  "
$ws.Range("D3").Value = "2.538.25"
$ws.Range("E3").Value = "  -3.28%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'298.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").Value = "'93.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.79%  "
$ws.Range("E7").Value = "  -3.71%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -5.16%  "
$ws.Range("D10").Value = "'35.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.11%  "
$ws.Range("D11").Value = "'0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.40%  "
$ws.Range("D12").Value = "'7.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.64%  "
$ws.Range("E13").Value = "  +4.66%  "
$ws.Range("D14").Value = "2.925.74"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("D15").Value = "2.557.69"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("E16").Value = "  -5.08%  "
$ws.Range("D17").Value = "'14.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.93%  "
$ws.Range("D18").Value = "42.720.60"
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("D19").Value = "'12.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").Value = "  -3.62%  "
$ws.Range("E21").Value = "  -3.26%  "
$ws.Range("E22").Value = "  -4.40%  "
$ws.Range("D23").Value = "'255.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -9.50%  "
$ws.Range("D24").Value = "'2.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.36%  "
$ws.Range("D25").Value = "'29.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("E26").Value = "  -6.60%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'36.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").Value = "'5.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("D32").Value = "'152.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("E33").Value = "  -6.23%  "
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").Value = "'3.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.29%  "
$ws.Range("D36").Value = "'0.0791"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.79%  "
$ws.Range("E37").Value = "  -7.61%  "
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").Value = "'24.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.26%  "
$ws.Range("D40").Value = "'16.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.53%  "
$ws.Range("E41").Value = "  -5.80%  "
$ws.Range("E42").Value = "  -4.15%  "
$ws.Range("E43").Value = "  -4.89%  "
$ws.Range("D44").Value = "2.077.24"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'84.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.13%  "
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").Value = "2.783.14"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").Value = "'104.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("E51").Value = "  -4.95%  "
